$p = $ppt.ActivePresentation

# Map of slide index -> chart shape index -> titles to apply.
# Chart shape is always the one whose .HasChart is true on that slide.
#   slide 7  : chart1 -> "EdTech Market Size Comparison" (axis titles already authored, leave as-is)
#   slide 10 : chart2 -> "AI Cost-Benefit Analysis" / Year / SGD (Millions)
#   slide 12 : chart3 -> "VR/AR Cost-Benefit Analysis" / Year / SGD (Millions)
#   slide 14 : chart4 -> "5G & Edge Computing Cost-Benefit Analysis" / Year / SGD (Millions)
#   slide 16 : chart5 -> "Blockchain Cost-Benefit Analysis" / Year / SGD (Millions)

function Get-ChartShape($slide) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasChart) {
            return $shp.Chart
        }
    }
    return $null
}

# --- Slide 7: EdTech Market Size Comparison ---------------------------------
$slide = $p.Slides.Item(7)
$chart = Get-ChartShape $slide
$chart.HasTitle = $true
$chart.ChartTitle.Text = "EdTech Market Size Comparison"

# --- Slide 10: AI Cost-Benefit Analysis --------------------------------------
$slide = $p.Slides.Item(10)
$chart = Get-ChartShape $slide
$chart.HasTitle = $true
$chart.ChartTitle.Text = "AI Cost-Benefit Analysis"

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Year"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "SGD (Millions)"

# --- Slide 12: VR/AR Cost-Benefit Analysis -----------------------------------
$slide = $p.Slides.Item(12)
$chart = Get-ChartShape $slide
$chart.HasTitle = $true
$chart.ChartTitle.Text = "VR/AR Cost-Benefit Analysis"

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Year"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "SGD (Millions)"

# --- Slide 14: 5G & Edge Computing Cost-Benefit Analysis ---------------------
$slide = $p.Slides.Item(14)
$chart = Get-ChartShape $slide
$chart.HasTitle = $true
$chart.ChartTitle.Text = "5G & Edge Computing Cost-Benefit Analysis"

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Year"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "SGD (Millions)"

# --- Slide 16: Blockchain Cost-Benefit Analysis ------------------------------
$slide = $p.Slides.Item(16)
$chart = Get-ChartShape $slide
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Blockchain Cost-Benefit Analysis"

$catAxis = $chart.Axes(1)
$catAxis.HasTitle = $true
$catAxis.AxisTitle.Text = "Year"

$valAxis = $chart.Axes(2)
$valAxis.HasTitle = $true
$valAxis.AxisTitle.Text = "SGD (Millions)"
